$d = $word.ActiveDocument

# Locate the "forest" / "and " run boundary (the two runs are adjacent with
# no space between them: "...forestand basins...").
$locator = $d.Content
$found = $locator.Find.Execute("forestand ")

if ($found) {
    # Narrow the range down to just the "and " run (last 4 characters of the
    # match), so the replace only touches that run and doesn't merge it with
    # the bold "forest" run.
    $andStart = $locator.End - 4
    $andEnd = $locator.End
    $r = $d.Range($andStart, $andEnd)
    $r.Find.Execute("and ", $false, $false, $false, $false, $false, $true, 1, $false, " and ", 2)
}
